# Lab6 report & little updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C1: "Epsilon " -> "Episode"
$ws.Range("C1").Value = "Episode"

# New cell F2 = 10000 (Capacity column was previously blank for the DQN row)
$ws.Range("F2").Value = 10000

# New row 4 - a third experiment run (DQN / Adam), mirrors row 2's formatting
$ws.Range("A2:M2").Copy()
$ws.Range("A4:M4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(4).RowHeight = 21.75

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "DQN"
$ws.Range("C4").Value = 2000
$ws.Range("D4").Value = 0.0005
$ws.Range("E4").Value = 128
$ws.Range("F4").Value = 10000
$ws.Range("G4").Value = "Adam"
$ws.Range("H4").Value = 0.995
$ws.Range("I4").Value = 0.01
$ws.Range("J4").Value = 0.99
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 0.001

# Update selection to match final saved state
$ws.Range("N6").Select()
